{"js": "// Replace each two-digit multiplication expression with its updated value.\n// Each old formula occurs exactly once in the document, so a body-wide\n// search-and-replace keyed on the distinctive old text is unambiguous.\nconst replacements = [\n  [\"44\u00d767=2948\", \"26\u00d715=390\"],\n  [\"42\u00d732=1344\", \"91\u00d727=2457\"],\n  [\"35\u00d713=455\", \"39\u00d729=1131\"],\n  [\"43\u00d750=2150\", \"25\u00d733=825\"],\n  [\"68\u00d755=3740\", \"15\u00d751=765\"],\n  [\"79\u00d781=6399\", \"32\u00d797=3104\"],\n  [\"95\u00d714=1330\", \"48\u00d782=3936\"],\n  [\"50\u00d770=3500\", \"18\u00d774=1332\"],\n  [\"40\u00d768=2720\", \"31\u00d758=1798\"],\n  [\"77\u00d749=3773\", \"79\u00d726=2054\"],\n  [\"89\u00d760=5340\", \"74\u00d749=3626\"],\n  [\"78\u00d798=7644\", \"74\u00d734=2516\"],\n  [\"16\u00d750=800\", \"23\u00d731=713\"],\n  [\"90\u00d759=5310\", \"69\u00d734=2346\"],\n  [\"29\u00d744=1276\", \"40\u00d756=2240\"],\n  [\"12\u00d745=540\", \"72\u00d763=4536\"],\n  [\"29\u00d747=1363\", \"15\u00d797=1455\"],\n  [\"42\u00d733=1386\", \"97\u00d724=2328\"],\n  [\"37\u00d742=1554\", \"72\u00d723=1656\"],\n  [\"74\u00d747=3478\", \"13\u00d791=1183\"],\n  [\"22\u00d768=1496\", \"86\u00d730=2580\"],\n  [\"55\u00d746=2530\", \"92\u00d776=6992\"],\n  [\"39\u00d739=1521\", \"55\u00d750=2750\"],\n  [\"62\u00d768=4216\", \"25\u00d771=1775\"],\n  [\"39\u00d749=1911\", \"95\u00d749=4655\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each two-digit multiplication expression with its updated value.\n# Each old formula occurs exactly once in the document, so Find/Replace keyed\n# on the distinctive old text is unambiguous. wdFindContinue = 1, wdReplaceOne = 2.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"44\u00d767=2948\", \"26\u00d715=390\"),\n    @(\"42\u00d732=1344\", \"91\u00d727=2457\"),\n    @(\"35\u00d713=455\", \"39\u00d729=1131\"),\n    @(\"43\u00d750=2150\", \"25\u00d733=825\"),\n    @(\"68\u00d755=3740\", \"15\u00d751=765\"),\n    @(\"79\u00d781=6399\", \"32\u00d797=3104\"),\n    @(\"95\u00d714=1330\", \"48\u00d782=3936\"),\n    @(\"50\u00d770=3500\", \"18\u00d774=1332\"),\n    @(\"40\u00d768=2720\", \"31\u00d758=1798\"),\n    @(\"77\u00d749=3773\", \"79\u00d726=2054\"),\n    @(\"89\u00d760=5340\", \"74\u00d749=3626\"),\n    @(\"78\u00d798=7644\", \"74\u00d734=2516\"),\n    @(\"16\u00d750=800\", \"23\u00d731=713\"),\n    @(\"90\u00d759=5310\", \"69\u00d734=2346\"),\n    @(\"29\u00d744=1276\", \"40\u00d756=2240\"),\n    @(\"12\u00d745=540\", \"72\u00d763=4536\"),\n    @(\"29\u00d747=1363\", \"15\u00d797=1455\"),\n    @(\"42\u00d733=1386\", \"97\u00d724=2328\"),\n    @(\"37\u00d742=1554\", \"72\u00d723=1656\"),\n    @(\"74\u00d747=3478\", \"13\u00d791=1183\"),\n    @(\"22\u00d768=1496\", \"86\u00d730=2580\"),\n    @(\"55\u00d746=2530\", \"92\u00d776=6992\"),\n    @(\"39\u00d739=1521\", \"55\u00d750=2750\"),\n    @(\"62\u00d768=4216\", \"25\u00d771=1775\"),\n    @(\"39\u00d749=1911\", \"95\u00d749=4655\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $ok) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
